$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.806.58"
$ws.Range("E2").Value = "  +4.78%  "
$ws.Range("D3").Value = "2.351.62"
$ws.Range("E3").Value = "  +4.49%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").Value = "'307.26"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'98.58"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("D10").Value = "'36.02"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("E12").Value = "  +3.35%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "2.708.20"
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("D15").Value = "2.359.90"
$ws.Range("E15").Value = "  +5.11%  "
$ws.Range("E16").Value = "  +5.25%  "
$ws.Range("D17").Value = "'0.834"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "46.702.09"
$ws.Range("E18").Value = "  +5.08%  "
$ws.Range("D19").Value = "'13.75"
$ws.Range("E19").Value = "  +16.83%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'66.83"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "'245.18"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'42.02"
$ws.Range("E27").Value = "  +12.96%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'9.94"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").Value = "'20.28"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'5.79"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").Value = "'152.50"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").Value = "'0.0820"
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'1.85"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'4.09"
$ws.Range("E39").Value = "  +8.14%  "
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").Value = "'14.00"
$ws.Range("E42").Value = "  -8.21%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'1.96"
$ws.Range("E44").Value = "  +9.97%  "
$ws.Range("D45").Value = "1.831.41"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "'0.199"
$ws.Range("E46").Value = "  +6.07%  "
$ws.Range("D47").Value = "'81.23"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "'73.67"
$ws.Range("E48").Value = "  +7.20%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'99.08"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'55.61"
$ws.Range("E51").Value = "  +2.82%  "
